# map_GTAP_format.xlsx - "reformat indexes" edit
#
# On the "Cost structure" sheet, the CI and Tax blocks (column A) are split
# into separate _imp (import) and _dom (domestic) sub-categories:
#   A1:A2  "CI"  -> "CI_imp"
#   A3:A4  "CI"  -> "CI_dom"
#   A28:A29 "Tax" -> "Tax_imp"
#   A30:A31 "Tax" -> "Tax_dom"
# Column B (and every other row/sheet) keeps its original text.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost structure")

$ws.Range("A1").Value = "CI_imp"
$ws.Range("A2").Value = "CI_imp"
$ws.Range("A3").Value = "CI_dom"
$ws.Range("A4").Value = "CI_dom"

$ws.Range("A28").Value = "Tax_imp"
$ws.Range("A29").Value = "Tax_imp"
$ws.Range("A30").Value = "Tax_dom"
$ws.Range("A31").Value = "Tax_dom"

# The author also navigated/selected a different cell before saving
# (selection moved from B46 to A31).
$ws.Activate()
$ws.Range("A31").Select()

$wb.Save()
